# Completed SM table of all socioecon models
#
# Adds the "Province level (categorical)" section (rows 41-51) and its
# corresponding "Final candidate model set" (rows 52-60) to the
# "Table S21" worksheet, mirroring the existing commune-level section
# (rows 1-40).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Table S21")

# --- Copy the existing cell formatting onto the new rows first, so the
#     new cells reuse the same styles already present in the sheet -------

# A41 mirrors the top-level section header style used by A1 ("Commune-level ").
$ws.Range("A1").Copy() | Out-Null
$ws.Range("A41").PasteSpecial(-4122) | Out-Null

# A42:A51 mirror the row-label style used for the "set" rows (e.g. A12).
$ws.Range("A12").Copy() | Out-Null
$ws.Range("A42:A51").PasteSpecial(-4122) | Out-Null

# B42:B51 mirror the style used in the "Final candidate model set" B column (e.g. B31).
$ws.Range("B31").Copy() | Out-Null
$ws.Range("B42:B51").PasteSpecial(-4122) | Out-Null

# A52:A60 mirror the same row-label style as A42:A51 above.
$ws.Range("A12").Copy() | Out-Null
$ws.Range("A52:A60").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# --- New section header (mirrors A1 "Commune-level ") ---
$ws.Range("A41").Value = "Province level (categorical)"

# --- Set-specific rows for the categorical (province-level) models ---
$ws.Range("A42").Value = "popden.mcat"
$ws.Range("B42").Value = "Population density"

$ws.Range("A43").Value = "socjus.mcat"
$ws.Range("B43").Value = "Land conflicts"
$ws.Range("C43").Value = "Criminal cases"

$ws.Range("A44").Value = "mig.mcat"
$ws.Range("B44").Value = "In-migration"
$ws.Range("C44").Value = "Out-migration"

$ws.Range("A45").Value = "edu.mcat"
$ws.Range("B45").Value = "Males in school"

$ws.Range("A46").Value = "emp.mcat"
$ws.Range("B46").Value = "Primary sector"
$ws.Range("C46").Value = "Secondary sector"

$ws.Range("A47").Value = "econ.mcat"
$ws.Range("B47").Value = "No farmland"
$ws.Range("C47").Value = "Owns pigs"

$ws.Range("A48").Value = "acc.mcat"
$ws.Range("B48").Value = "Distance to school"

$ws.Range("A49").Value = "elev.mcat"
$ws.Range("B49").Value = "Elevation"

$ws.Range("A50").Value = "hum.mcat"
$ws.Range("B50").Value = "Distance to Int'l border"
$ws.Range("C50").Value = "Distance to provincial capital"

$ws.Range("A51").Value = "areas.mcat"
$ws.Range("B51").Value = "Presence of ELC"
$ws.Range("C51").Value = "Presence of PAs"

# --- Final candidate model set (province-level, categorical) ---
$ws.Range("A52").Value = "m1"
$ws.Range("B52").Value = "Population density"
$ws.Range("C52").Value = "Elevation"
$ws.Range("D52").Value = "Distance to Int'l border"
$ws.Range("E52").Value = "Distance to provincial capital"
$ws.Range("F52").Value = "Presence of ELC"
$ws.Range("G52").Value = "Presence of PA"

$ws.Range("A53").Value = "m2"
$ws.Range("B53").Value = "Males in school"
$ws.Range("C53").Value = "Elevation"
$ws.Range("D53").Value = "Distance to Int'l border"
$ws.Range("E53").Value = "Distance to provincial capital"
$ws.Range("F53").Value = "Presence of ELC"
$ws.Range("G53").Value = "Presence of PA"

$ws.Range("A54").Value = "m3"
$ws.Range("B54").Value = "Primary sector"
$ws.Range("C54").Value = "Elevation"
$ws.Range("D54").Value = "Distance to Int'l border"
$ws.Range("E54").Value = "Distance to provincial capital"
$ws.Range("F54").Value = "Presence of ELC"
$ws.Range("G54").Value = "Presence of PA"

$ws.Range("A55").Value = "m4"
$ws.Range("B55").Value = "Pig ownership"
$ws.Range("C55").Value = "Elevation"
$ws.Range("D55").Value = "Distance to Int'l border"
$ws.Range("E55").Value = "Distance to provincial capital"
$ws.Range("F55").Value = "Presence of ELC"
$ws.Range("G55").Value = "Presence of PA"

$ws.Range("A56").Value = "m5"
$ws.Range("B56").Value = "Distance to school"
$ws.Range("C56").Value = "Elevation"
$ws.Range("D56").Value = "Distance to Int'l border"
$ws.Range("E56").Value = "Distance to provincial capital"
$ws.Range("F56").Value = "Presence of ELC"
$ws.Range("G56").Value = "Presence of PA"

$ws.Range("A57").Value = "m6"
$ws.Range("B57").Value = "Criminal cases"
$ws.Range("C57").Value = "Elevation"
$ws.Range("D57").Value = "Distance to Int'l border"
$ws.Range("E57").Value = "Distance to provincial capital"
$ws.Range("F57").Value = "Presence of ELC"
$ws.Range("G57").Value = "Presence of PA"

$ws.Range("A58").Value = "m7"
$ws.Range("B58").Value = "Out-migration"
$ws.Range("C58").Value = "Elevation"
$ws.Range("D58").Value = "Distance to Int'l border"
$ws.Range("E58").Value = "Distance to provincial capital"
$ws.Range("F58").Value = "Presence of ELC"
$ws.Range("G58").Value = "Presence of PA"

$ws.Range("A59").Value = "m8"
$ws.Range("B59").Value = "Males in school"
$ws.Range("C59").Value = "distance to school"
$ws.Range("D59").Value = "Elevation"
$ws.Range("E59").Value = "Distance to Int'l border"
$ws.Range("F59").Value = "Distance to provincial capital"
$ws.Range("G59").Value = "Presence of ELC"
$ws.Range("H59").Value = "Presence of PA"

$ws.Range("A60").Value = "m9"
$ws.Range("B60").Value = "Primary sector"
$ws.Range("C60").Value = "Out-migration"
$ws.Range("D60").Value = "Elevation"
$ws.Range("E60").Value = "Distance to Int'l border"
$ws.Range("F60").Value = "Distance to provincial capital"
$ws.Range("G60").Value = "Presence of ELC"
$ws.Range("H60").Value = "Presence of PA"

# --- Match the authored view state: the header row scrolled back to the
#     top and the new D51 cell selected. ---
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("D51").Select()
